# Add new "marine-first" variant rows to the habitat translation table.
# Each new row is inserted directly below an existing row that shares the
# same translated (column B) value, just with the terms reordered so that
# "marine" comes first in the original term (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert the five new (still empty) rows. We insert bottom-to-top
# (using the original, pre-insert row numbers) so that each insertion point
# is unaffected by the ones that still need to happen further up the sheet.
$ws.Rows.Item(16).Insert()   # -> new row 16, ends up as final row 20
$ws.Rows.Item(14).Insert()   # -> new row 14, ends up as final row 17
$ws.Rows.Item(10).Insert()   # -> new row 10, ends up as final row 12
$ws.Rows.Item(9).Insert()    # -> new row 9,  ends up as final row 10
$ws.Rows.Item(8).Insert()    # -> new row 8,  ends up as final row 8

# Step 2: fill in the values, in the same order the new strings were
# originally authored, so the shared-string table ends up in that order.
$ws.Range("A20").Value = "marine; terrestrial"
$ws.Range("B20").Value = "terrestrial; marine"

$ws.Range("A12").Value = "marine; terrestrial; freshwater"
$ws.Range("B12").Value = "terrestrial; freshwater; marine"

$ws.Range("A17").Value = "marine; freshwater"
$ws.Range("B17").Value = "freshwater; marine"

$ws.Range("A10").Value = "marine; brackish"
$ws.Range("B10").Value = "brackish; marine"

$ws.Range("A8").Value = "marine; freshwater; brackish"
$ws.Range("B8").Value = "freshwater; brackish; marine"

$ws.Range("A14").Select()
